# Add a new weekly price record for Vega Monumental Concepción - Espinaca.
# This inserts a new row at position 105 (shifting existing rows 105-125 down
# to 106-126) and populates it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(105).Insert()

$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 45093
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = 100112012
$ws.Cells.Item(105, 7).Value = "Espinaca"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 40
$ws.Cells.Item(105, 11).Value = 6500
$ws.Cells.Item(105, 12).Value = 7000
$ws.Cells.Item(105, 13).Value = 6750
$ws.Cells.Item(105, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(105, 15).Value = "Región Metropolitana"
$ws.Cells.Item(105, 16).Value = 675
$ws.Cells.Item(105, 17).Value = 10
$ws.Cells.Item(105, 18).Value = "Hortaliza"
